# Update Developer guide images
#  - bump the cached datetimeFigureOut field text from 11/11/2018 to 11/12/2018
#    on the slide master and every slide layout
#  - rename AddressBookChangedEvent -> HealthBookChangedEvent (and the
#    matching handler name) in the sequence-diagram textboxes on slide 1

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder fields (slide master + all 11 custom layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "11/11/2018") {
                $tr.Text = "11/12/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}

# ---------------------------------------------------------------------
# 2) Rename AddressBookChangedEvent -> HealthBookChangedEvent on slide 1
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $tr = $sh.TextFrame.TextRange
    $full = $tr.Text

    $oldEvt = "AddressBookChangedEvent"
    $newEvt = "HealthBookChangedEvent"
    $oldHandler = "handleAddresssBookChangedEvent"
    $newHandler = "handleHealthBookChangedEvent"

    if ($full -eq "post(AddressBookChangedEvent)") {
        $pos = $full.IndexOf($oldEvt) + 1
        $sub = $tr.Characters($pos, $oldEvt.Length)
        $sub.Text = $newEvt
    }
    elseif ($full -eq "handleAddresssBookChangedEvent()") {
        $pos = $full.IndexOf($oldHandler) + 1
        $sub = $tr.Characters($pos, $oldHandler.Length)
        $sub.Text = $newHandler
    }
}
